$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -11.9585
$ws.Range("C6").Value = -13.0543
$ws.Range("C7").Value = -12.29349999999999
$ws.Range("D7").Value = -7.266699999999999
$ws.Range("D12").Value = -7.152999999999997
$ws.Range("D15").Value = -8.560599999999996
$ws.Range("C16").Value = -14.07889999999999
$ws.Range("C20").Value = -12.65290000000001
$ws.Range("D20").Value = -7.994899999999993
$ws.Range("D21").Value = -7.970899999999998
$ws.Range("D22").Value = -8.145500000000004
$ws.Range("D23").Value = -7.393299999999998
$ws.Range("C28").Value = -12.0952
$ws.Range("C29").Value = -11.11510000000001
$ws.Range("D29").Value = -7.1067
$ws.Range("C32").Value = -12.8587
$ws.Range("D34").Value = -7.885400000000002
$ws.Range("C40").Value = -11.69630000000001
$ws.Range("D42").Value = -8.494600000000004
$ws.Range("D43").Value = -8.259700000000002
$ws.Range("D44").Value = -7.080800000000003
$ws.Range("D45").Value = -7.532799999999997
$ws.Range("C46").Value = -14.22319999999999
$ws.Range("D46").Value = -8.351300000000002
$ws.Range("D50").Value = -7.998499999999998
$ws.Range("C51").Value = -11.7747
$ws.Range("D51").Value = -7.510599999999998
$ws.Range("C52").Value = -11.1393
$ws.Range("C57").Value = -14.31849999999999
$ws.Range("C59").Value = -12.6794
$ws.Range("C62").Value = -13.90069999999999
$ws.Range("C66").Value = -11.31250000000001
$ws.Range("D66").Value = -7.415199999999999
$ws.Range("D67").Value = -6.4072
$ws.Range("C73").Value = -10.99090000000001
$ws.Range("C74").Value = -11.91390000000001
$ws.Range("D79").Value = -6.381600000000002
$ws.Range("D84").Value = -8.7896
$ws.Range("C92").Value = -11.2721
$ws.Range("D92").Value = -6.559500000000003
$ws.Range("D97").Value = -8.237400000000003
$ws.Range("C100").Value = -12.2921
